$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.727.23"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.648.12"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.96"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.14"
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "1.877.39"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "1.638.00"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.16"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "27.632.50"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.33"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.10"
$ws.Range("E23").Value = "  +10.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.14"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.98"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.63"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "1.434.94"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.574"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.896"
$ws.Range("E40").Value = "  +14.19%  "
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.43"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "1.786.20"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.69"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.48"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0990"
$ws.Range("E51").Value = "  -1.98%  "
